$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 2).Value = 16.52184105490078
$ws.Cells.Item(2, 4).Value = 7.877306581233581
$ws.Cells.Item(2, 5).Value = 25.59877912509573
$ws.Cells.Item(2, 6).Value = 45.47921598967443
$ws.Cells.Item(2, 7).Value = 3.645392153895943
$ws.Cells.Item(2, 9).Value = 25.93575466912673
$ws.Cells.Item(2, 12).Value = 12.28798852929961
$ws.Cells.Item(2, 13).Value = 17.1390882066967
$ws.Cells.Item(2, 14).Value = 19.04854897325594

$ws.Cells.Item(3, 2).Value = 16.24611252694594
$ws.Cells.Item(3, 4).Value = 7.786876854999162
$ws.Cells.Item(3, 5).Value = 24.0031485881934
$ws.Cells.Item(3, 6).Value = 44.24575825043235
$ws.Cells.Item(3, 7).Value = 3.653731593812685
$ws.Cells.Item(3, 9).Value = 26.10612614482224
$ws.Cells.Item(3, 12).Value = 12.11960071415472
$ws.Cells.Item(3, 13).Value = 16.99272164858797
$ws.Cells.Item(3, 14).Value = 19.12881201816397

$ws.Cells.Item(4, 2).Value = 16.07848944164112
$ws.Cells.Item(4, 4).Value = 7.73451838414412
$ws.Cells.Item(4, 5).Value = 22.96773923048344
$ws.Cells.Item(4, 6).Value = 43.49324691037108
$ws.Cells.Item(4, 7).Value = 3.659092497241657
$ws.Cells.Item(4, 9).Value = 26.21667587669087
$ws.Cells.Item(4, 12).Value = 12.01893127232456
$ws.Cells.Item(4, 13).Value = 16.9065511017785
$ws.Cells.Item(4, 14).Value = 19.18020693071199

$ws.Cells.Item(5, 2).Value = 16.01068919797066
$ws.Cells.Item(5, 4).Value = 7.713993559057991
$ws.Cells.Item(5, 5).Value = 22.53182580684014
$ws.Cells.Item(5, 6).Value = 43.18830315651833
$ws.Cells.Item(5, 7).Value = 3.661337985014467
$ws.Cells.Item(5, 9).Value = 26.26321730582593
$ws.Cells.Item(5, 12).Value = 11.97863527749824
$ws.Cells.Item(5, 13).Value = 16.87239474255985
$ws.Cells.Item(5, 14).Value = 19.20168549587521

$ws.Cells.Item(6, 2).Value = 15.99946418861927
$ws.Cells.Item(6, 4).Value = 7.71063486237626
$ws.Cells.Item(6, 5).Value = 22.45859979035257
$ws.Cells.Item(6, 6).Value = 43.13778548085344
$ws.Cells.Item(6, 7).Value = 3.661714535314801
$ws.Cells.Item(6, 9).Value = 26.27103549353134
$ws.Cells.Item(6, 12).Value = 11.97198933777373
$ws.Cells.Item(6, 13).Value = 16.86678181812932
$ws.Cells.Item(6, 14).Value = 19.20528438947793

$ws.Cells.Item(7, 2).Value = 16.07757289506423
$ws.Cells.Item(7, 4).Value = 7.734238273754002
$ws.Cells.Item(7, 5).Value = 22.96191688729791
$ws.Cells.Item(7, 6).Value = 43.48912674330745
$ws.Cells.Item(7, 7).Value = 3.659122533645292
$ws.Cells.Item(7, 9).Value = 26.21729751450737
$ws.Cells.Item(7, 12).Value = 12.01838482431503
$ws.Cells.Item(7, 13).Value = 16.90608653833797
$ws.Cells.Item(7, 14).Value = 19.18049442861843

$ws.Cells.Item(8, 2).Value = 16.42647388089823
$ws.Cells.Item(8, 4).Value = 7.84547941573805
$ws.Cells.Item(8, 5).Value = 25.06017358445896
$ws.Cells.Item(8, 6).Value = 45.05320838257759
$ws.Cells.Item(8, 7).Value = 3.648217936430425
$ws.Cells.Item(8, 9).Value = 25.99326421499789
$ws.Cells.Item(8, 12).Value = 12.22939004874122
$ws.Cells.Item(8, 13).Value = 17.08787100649804
$ws.Cells.Item(8, 14).Value = 19.07578714319735

$ws.Cells.Item(9, 2).Value = 17.12018753189969
$ws.Cells.Item(9, 4).Value = 8.087951527670935
$ws.Cells.Item(9, 5).Value = 28.7327542504498
$ws.Cells.Item(9, 6).Value = 48.13783410247628
$ws.Cells.Item(9, 7).Value = 3.62872300301114
$ws.Cells.Item(9, 9).Value = 25.60118835886752
$ws.Cells.Item(9, 12).Value = 12.66279461513936
$ws.Cells.Item(9, 13).Value = 17.47236101918654
$ws.Cells.Item(9, 14).Value = 18.88706749896543

$ws.Cells.Item(10, 2).Value = 17.63068892543015
$ws.Cells.Item(10, 4).Value = 8.279743900927208
$ws.Cells.Item(10, 5).Value = 31.16220475206405
$ws.Cells.Item(10, 6).Value = 50.38832148434742
$ws.Cells.Item(10, 7).Value = 3.615524868415223
$ws.Cells.Item(10, 9).Value = 25.3421053732493
$ws.Cells.Item(10, 12).Value = 12.99048701702924
$ws.Cells.Item(10, 13).Value = 17.77003133685057
$ws.Cells.Item(10, 14).Value = 18.75832179302512

$ws.Cells.Item(11, 2).Value = 17.86209472344641
$ws.Cells.Item(11, 4).Value = 8.369677581826339
$ws.Cells.Item(11, 5).Value = 32.20943051758615
$ws.Cells.Item(11, 6).Value = 51.40364449528278
$ws.Cells.Item(11, 7).Value = 3.609759066165345
$ws.Cells.Item(11, 9).Value = 25.23057633466069
$ws.Cells.Item(11, 12).Value = 13.14099865863129
$ws.Cells.Item(11, 13).Value = 17.90832128551209
$ws.Cells.Item(11, 14).Value = 18.70185519763612

$ws.Cells.Item(12, 2).Value = 17.94952272247748
$ws.Cells.Item(12, 4).Value = 8.404094384505644
$ws.Cells.Item(12, 5).Value = 32.59769339871195
$ws.Cells.Item(12, 6).Value = 51.78652823368607
$ws.Cells.Item(12, 7).Value = 3.607609471538992
$ws.Cells.Item(12, 9).Value = 25.18925768576582
$ws.Cells.Item(12, 12).Value = 13.1981531502609
$ws.Cells.Item(12, 13).Value = 17.96106507487339
$ws.Cells.Item(12, 14).Value = 18.68077100878137

$ws.Cells.Item(13, 2).Value = 17.93070371583143
$ws.Cells.Item(13, 4).Value = 8.396666470866711
$ws.Cells.Item(13, 5).Value = 32.51444253427238
$ws.Cells.Item(13, 6).Value = 51.70414416816115
$ws.Cells.Item(13, 7).Value = 3.608070929620744
$ws.Cells.Item(13, 9).Value = 25.19811562058655
$ws.Cells.Item(13, 12).Value = 13.18583759766993
$ws.Cells.Item(13, 13).Value = 17.94968961710582
$ws.Cells.Item(13, 14).Value = 18.6852986480308

$ws.Cells.Item(14, 2).Value = 17.86929194955412
$ws.Cells.Item(14, 4).Value = 8.372501988034516
$ws.Cells.Item(14, 5).Value = 32.24153924518018
$ws.Cells.Item(14, 6).Value = 51.43517866237396
$ws.Cells.Item(14, 7).Value = 3.609581543125969
$ws.Cells.Item(14, 9).Value = 25.22715865784229
$ws.Cells.Item(14, 12).Value = 13.14569788069313
$ws.Cells.Item(14, 13).Value = 17.91265322158731
$ws.Cells.Item(14, 14).Value = 18.70011462720629

$ws.Cells.Item(15, 2).Value = 17.83164700665748
$ws.Cells.Item(15, 4).Value = 8.35774677654353
$ws.Cells.Item(15, 5).Value = 32.07329845063423
$ws.Cells.Item(15, 6).Value = 51.27021069707469
$ws.Cells.Item(15, 7).Value = 3.610511224479227
$ws.Cells.Item(15, 9).Value = 25.24506766479066
$ws.Cells.Item(15, 12).Value = 13.12113042915404
$ws.Cells.Item(15, 13).Value = 17.89001525718927
$ws.Cells.Item(15, 14).Value = 18.70922860744127

$ws.Cells.Item(16, 2).Value = 17.61554235608087
$ws.Cells.Item(16, 4).Value = 8.273918419843888
$ws.Cells.Item(16, 5).Value = 31.09260282692163
$ws.Cells.Item(16, 6).Value = 50.3217644279095
$ws.Cells.Item(16, 7).Value = 3.615906430084027
$ws.Cells.Item(16, 9).Value = 25.34952182242791
$ws.Cells.Item(16, 12).Value = 12.98067559557542
$ws.Cells.Item(16, 13).Value = 17.76104848482242
$ws.Cells.Item(16, 14).Value = 18.7620539813339

$ws.Cells.Item(17, 2).Value = 17.48269880732849
$ws.Cells.Item(17, 4).Value = 8.223162588847172
$ws.Cells.Item(17, 5).Value = 30.47615884908106
$ws.Cells.Item(17, 6).Value = 49.73747938576305
$ws.Cells.Item(17, 7).Value = 3.619276879207482
$ws.Cells.Item(17, 9).Value = 25.41522543438906
$ws.Cells.Item(17, 12).Value = 12.89484587119296
$ws.Cells.Item(17, 13).Value = 17.68264219238977
$ws.Cells.Item(17, 14).Value = 18.79499609851954

$ws.Cells.Item(18, 2).Value = 17.40621786196076
$ws.Cells.Item(18, 4).Value = 8.194223601949975
$ws.Cells.Item(18, 5).Value = 30.11614022147575
$ws.Cells.Item(18, 6).Value = 49.40064108351774
$ws.Cells.Item(18, 7).Value = 3.621237910298508
$ws.Cells.Item(18, 9).Value = 25.45361182557333
$ws.Cells.Item(18, 12).Value = 12.84561753656557
$ws.Cells.Item(18, 13).Value = 17.63781815900192
$ws.Cells.Item(18, 14).Value = 18.81414145475249

$ws.Cells.Item(19, 2).Value = 17.38031276695196
$ws.Cells.Item(19, 4).Value = 8.184469820677213
$ws.Cells.Item(19, 5).Value = 29.99330605227731
$ws.Cells.Item(19, 6).Value = 49.28647289390551
$ws.Cells.Item(19, 7).Value = 3.621905749002637
$ws.Cells.Item(19, 9).Value = 25.4667109189071
$ws.Cells.Item(19, 12).Value = 12.82897497288759
$ws.Cells.Item(19, 13).Value = 17.62268956796439
$ws.Cells.Item(19, 14).Value = 18.820657844437

$ws.Cells.Item(20, 2).Value = 17.49684837369021
$ws.Cells.Item(20, 4).Value = 8.228539484088845
$ws.Cells.Item(20, 5).Value = 30.5423448526496
$ws.Cells.Item(20, 6).Value = 49.79976041428609
$ws.Cells.Item(20, 7).Value = 3.618915770132415
$ws.Cells.Item(20, 9).Value = 25.40816951599584
$ws.Cells.Item(20, 12).Value = 12.90396859578146
$ws.Cells.Item(20, 13).Value = 17.69096067205813
$ws.Cells.Item(20, 14).Value = 18.79146889542253

$ws.Cells.Item(21, 2).Value = 17.88733614670066
$ws.Cells.Item(21, 4).Value = 8.379590093792434
$ws.Cells.Item(21, 5).Value = 32.32192248672174
$ws.Cells.Item(21, 6).Value = 51.5142265941935
$ws.Cells.Item(21, 7).Value = 3.609136925889377
$ws.Cells.Item(21, 9).Value = 25.2186031408027
$ws.Cells.Item(21, 12).Value = 13.15748395483275
$ws.Cells.Item(21, 13).Value = 17.92352179432882
$ws.Cells.Item(21, 14).Value = 18.69575474273804

$ws.Cells.Item(22, 2).Value = 18.14133922581663
$ws.Cells.Item(22, 4).Value = 8.480401721610866
$ws.Cells.Item(22, 5).Value = 33.43662727072916
$ws.Cells.Item(22, 6).Value = 52.62526546674881
$ws.Cells.Item(22, 7).Value = 3.602942623923563
$ws.Cells.Item(22, 9).Value = 25.10004628916175
$ws.Cells.Item(22, 12).Value = 13.32407593715661
$ws.Cells.Item(22, 13).Value = 18.0776902052004
$ws.Cells.Item(22, 14).Value = 18.63493812870703

$ws.Cells.Item(23, 2).Value = 18.00590923013638
$ws.Cells.Item(23, 4).Value = 8.426413824809266
$ws.Cells.Item(23, 5).Value = 32.84609999044403
$ws.Cells.Item(23, 6).Value = 52.03326750665125
$ws.Cells.Item(23, 7).Value = 3.606230787943833
$ws.Cells.Item(23, 9).Value = 25.1628323817227
$ws.Cells.Item(23, 12).Value = 13.2350955967928
$ws.Cells.Item(23, 13).Value = 17.99522096471986
$ws.Cells.Item(23, 14).Value = 18.66723924974468

$ws.Cells.Item(24, 2).Value = 17.49045168307533
$ws.Cells.Item(24, 4).Value = 8.226107836491414
$ws.Cells.Item(24, 5).Value = 30.51243965789511
$ws.Cells.Item(24, 6).Value = 49.77160601310047
$ws.Cells.Item(24, 7).Value = 3.619078954890327
$ws.Cells.Item(24, 9).Value = 25.41135758923395
$ws.Cells.Item(24, 12).Value = 12.8998438460058
$ws.Cells.Item(24, 13).Value = 17.68719909718071
$ws.Cells.Item(24, 14).Value = 18.79306290049027

$ws.Cells.Item(25, 2).Value = 16.93202066270345
$ws.Cells.Item(25, 4).Value = 8.019867788802156
$ws.Cells.Item(25, 5).Value = 27.78652283400253
$ws.Cells.Item(25, 6).Value = 47.30431883265932
$ws.Cells.Item(25, 7).Value = 3.633797419199541
$ws.Cells.Item(25, 9).Value = 25.70218035090023
$ws.Cells.Item(25, 12).Value = 12.54372259737793
$ws.Cells.Item(25, 13).Value = 17.36553778724245
$ws.Cells.Item(25, 14).Value = 18.93636545723941
